$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 141 (pushes the existing 141:146 block down to 144:149,
# matching the trailing rows already present in the sheet).
$ws.Rows("141:143").Insert()

# Row 141: new weekly "Especial" quote (date 45147)
$ws.Range("A141").Value = 2
$ws.Range("B141").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C141").Value = "Coquimbo"
$ws.Range("D141").Value = 45147
$ws.Range("E141").Value = 4
$ws.Range("F141").Value = "Fruta"
$ws.Range("G141").Value = 100107
$ws.Range("H141").Value = "Otros"
$ws.Range("I141").Value = 100107002
$ws.Range("J141").Value = "Chirimoya"
$ws.Range("K141").Value = "Cultivar IV Región"
$ws.Range("L141").Value = "Especial"
$ws.Range("M141").Value = 200
$ws.Range("N141").Value = 22000
$ws.Range("O141").Value = 23000
$ws.Range("P141").Value = 22500
$ws.Range("Q141").Value = "`$/bandeja 10 kilos"
$ws.Range("R141").Value = "Provincia de Limarí"
$ws.Range("S141").Value = 2250
$ws.Range("T141").Value = 10

# Row 142: new weekly "Primera" quote (date 45147)
$ws.Range("A142").Value = 2
$ws.Range("B142").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C142").Value = "Coquimbo"
$ws.Range("D142").Value = 45147
$ws.Range("E142").Value = 4
$ws.Range("F142").Value = "Fruta"
$ws.Range("G142").Value = 100107
$ws.Range("H142").Value = "Otros"
$ws.Range("I142").Value = 100107002
$ws.Range("J142").Value = "Chirimoya"
$ws.Range("K142").Value = "Cultivar IV Región"
$ws.Range("L142").Value = "Primera"
$ws.Range("M142").Value = 300
$ws.Range("N142").Value = 19000
$ws.Range("O142").Value = 20000
$ws.Range("P142").Value = 19500
$ws.Range("Q142").Value = "`$/bandeja 10 kilos"
$ws.Range("R142").Value = "Provincia de Limarí"
$ws.Range("S142").Value = 1950
$ws.Range("T142").Value = 10

# Row 143: new weekly "Segunda" quote (date 45147)
$ws.Range("A143").Value = 2
$ws.Range("B143").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C143").Value = "Coquimbo"
$ws.Range("D143").Value = 45147
$ws.Range("E143").Value = 4
$ws.Range("F143").Value = "Fruta"
$ws.Range("G143").Value = 100107
$ws.Range("H143").Value = "Otros"
$ws.Range("I143").Value = 100107002
$ws.Range("J143").Value = "Chirimoya"
$ws.Range("K143").Value = "Cultivar IV Región"
$ws.Range("L143").Value = "Segunda"
$ws.Range("M143").Value = 240
$ws.Range("N143").Value = 16000
$ws.Range("O143").Value = 17000
$ws.Range("P143").Value = 16500
$ws.Range("Q143").Value = "`$/bandeja 10 kilos"
$ws.Range("R143").Value = "Provincia de Limarí"
$ws.Range("S143").Value = 1650
$ws.Range("T143").Value = 10
